$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new log rows (32 and 33) to the feed_logs sheet
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "2024-06-15 13:12:52"
$ws.Range("D32").Value = 200
$ws.Range("E32").Value = 7

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 2
$ws.Range("C33").Value = "2024-06-15 13:12:52"
$ws.Range("D33").Value = 200
$ws.Range("E33").Value = 0
